$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two removed data rows (delete higher row index first so the
# lower one's row number stays valid).
$ws.Rows.Item(28).Delete()   # "SC 92" row
$ws.Rows.Item(26).Delete()   # "RM 232" row

# Update the C-column "missing value" cells that changed between the two
# snapshots (some became known values, some became missing again).
$ws.Range("C3").Value = 11.2
$ws.Range("C5").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = ""
$ws.Range("C32").Value = 10.5
